# Rename the header labels on row 1 of each data table so that Power BI
# automatically treats the first row as the table header.
# Sheets 1, 2, 3 and 5 use "Ano <year>" labels (single year columns).
# Sheet 4 uses "Intervalo <period>" labels (year / interval columns).
# Sheet 6 only has a single year column (B1), also renamed to "Ano 2015".

$wb = $excel.ActiveWorkbook

$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $anoSheets) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($col in @("B", "C", "D", "E")) {
        $cell = $ws.Range($col + "1")
        $cell.Value = "Ano " + $cell.Value2
    }
}

$wsIntervalo = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
foreach ($col in @("B", "C", "D", "E")) {
    $cell = $wsIntervalo.Range($col + "1")
    $cell.Value = "Intervalo " + $cell.Value2
}

$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$custoCell = $wsCusto.Range("B1")
$custoCell.Value = "Ano " + $custoCell.Value2
